# Update the "汽車" (car) sheet:
#  - Insert a new "capacity" column after "name" (old B column), shifting the
#    existing owner/register_date/register_reason/acquire_value columns right
#  - Fix row 1 to hold proper header labels instead of duplicated row-2 data
#  - Append the standard trailing metadata columns (property_category,
#    category, date, legislator_name, legislator_id, source_file, index)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Capture the existing row 2 values before they get overwritten, since the
# columns need to shift right by one (new capacity column inserted at C).
$name          = $ws.Range("B2").Value()
$area          = $ws.Range("C2").Value()
$owner         = $ws.Range("D2").Value()
$registerDate  = $ws.Range("E2").Value()
$registerReason= $ws.Range("F2").Value()
$acquireValue  = $ws.Range("G2").Value()

# --- Row 1: proper header labels ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: shift existing values right, keep the same data ---
$ws.Range("B2").Value = $name
$ws.Range("C2").Value = $area
$ws.Range("D2").Value = $owner
$ws.Range("E2").Value = $registerDate
$ws.Range("F2").Value = $registerReason
$ws.Range("G2").Value = $acquireValue
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-12"
$ws.Range("K2").Value = "林鴻池"
$ws.Range("L2").Value = 1340
$ws.Range("M2").Value = "tmp67ea1"
$ws.Range("N2").Value = 67
